$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("M2").Value = 0.2303363333333333
$ws.Range("N2").Value = 0.691009
$ws.Range("O2").Value = 0.0420565315194687
$ws.Range("P2").Value = 0.0420565315194687
$ws.Range("Q2").Value = 0.08597910194011112
$ws.Range("R2").Value = 0.773811917461
$ws.Range("S2").Value = 0.01494798043935924
$ws.Range("T2").Value = 0.01494798043935924

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("O3").Value = 0.8440851393264226
$ws.Range("P3").Value = 0.8440851393264227
$ws.Range("Q3").Value = 1.725622147577333
$ws.Range("R3").Value = 15.530599328196
$ws.Range("S3").Value = 0.3000097177762838
$ws.Range("T3").Value = 0.3000097177762838

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3732763333333333
$ws.Range("H4").Value = 1.119829
$ws.Range("I4").Value = 0.3554258969843855
$ws.Range("J4").Value = 0.3554258969843855
$ws.Range("M4").Value = 0.6235823333333333
$ws.Range("N4").Value = 1.870747
$ws.Range("O4").Value = 0.1138583291541087
$ws.Range("P4").Value = 0.1138583291541087
$ws.Range("Q4").Value = 0.2327685269181111
$ws.Range("R4").Value = 2.094916742263
$ws.Range("S4").Value = 0.0404681987687425
$ws.Range("T4").Value = 0.0404681987687425

# Row 5
$ws.Range("G5").Value = 0.668317
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("M5").Value = 0.2303363333333333
$ws.Range("N5").Value = 0.691009
$ws.Range("O5").Value = 0.0420565315194687
$ws.Range("P5").Value = 0.0420565315194687
$ws.Range("Q5").Value = 0.1539376872843334
$ws.Range("R5").Value = 1.385439185559
$ws.Range("S5").Value = 0.02676298642906529
$ws.Range("T5").Value = 0.02676298642906529

# Row 6
$ws.Range("G6").Value = 0.668317
$ws.Range("I6").Value = 0.6363574327729865
$ws.Range("J6").Value = 0.6363574327729865
$ws.Range("O6").Value = 0.8440851393264226
$ws.Range("P6").Value = 0.8440851393264227
$ws.Range("S6").Value = 0.537139852303591
$ws.Range("T6").Value = 0.537139852303591

# Row 7
$ws.Range("G7").Value = 0.668317
$ws.Range("I7").Value = 0.6363574327729865
$ws.Range("J7").Value = 0.6363574327729865
$ws.Range("M7").Value = 0.6235823333333333
$ws.Range("N7").Value = 1.870747
$ws.Range("O7").Value = 0.1138583291541087
$ws.Range("P7").Value = 0.1138583291541087
$ws.Range("Q7").Value = 0.4167506742663334
$ws.Range("R7").Value = 3.750756068397
$ws.Range("S7").Value = 0.0724545940403303
$ws.Range("T7").Value = 0.0724545940403303

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008629333333333334
$ws.Range("H8").Value = 0.025888
$ws.Range("I8").Value = 0.008216670242627913
$ws.Range("J8").Value = 0.008216670242627911
$ws.Range("M8").Value = 0.2303363333333333
$ws.Range("N8").Value = 0.691009
$ws.Range("O8").Value = 0.0420565315194687
$ws.Range("P8").Value = 0.0420565315194687
$ws.Range("Q8").Value = 0.001987648999111111
$ws.Range("R8").Value = 0.017888840992
$ws.Range("S8").Value = 0.0003455646510441614
$ws.Range("T8").Value = 0.0003455646510441613

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008629333333333334
$ws.Range("H9").Value = 0.025888
$ws.Range("I9").Value = 0.008216670242627913
$ws.Range("J9").Value = 0.008216670242627911
$ws.Range("O9").Value = 0.8440851393264226
$ws.Range("P9").Value = 0.8440851393264227
$ws.Range("Q9").Value = 0.03989261410133334
$ws.Range("R9").Value = 0.359033526912
$ws.Range("S9").Value = 0.006935569246547852
$ws.Range("T9").Value = 0.006935569246547852

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008629333333333334
$ws.Range("H10").Value = 0.025888
$ws.Range("I10").Value = 0.008216670242627913
$ws.Range("J10").Value = 0.008216670242627911
$ws.Range("M10").Value = 0.6235823333333333
$ws.Range("N10").Value = 1.870747
$ws.Range("O10").Value = 0.1138583291541087
$ws.Range("P10").Value = 0.1138583291541087
$ws.Range("Q10").Value = 0.005381099815111112
$ws.Range("R10").Value = 0.048429898336
$ws.Range("S10").Value = 0.0009355363450358992
$ws.Range("T10").Value = 0.0009355363450358991
